$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 2: rename / add GPU timing columns ---
# (order matters for how shared strings get (re)used internally)
$ws.Range("F2").Value = "GPU (stream for mag/dir)"
$ws.Range("D2").Value = "GPU (stream for sobel x / y)"
$ws.Range("E2").Value = "changes"
$ws.Range("G2").Value = "changes"

# --- New data column F (GPU stream for mag/dir timings) ---
$ws.Range("F3").Value = 0.19964999999999999
$ws.Range("F4").Value = 0.19980999999999999
$ws.Range("F5").Value = 0.00083000000000000001
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 0.00109
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 0.0848
$ws.Range("F10").Value = 0.09206
$ws.Range("F11").Value = 0.11866

# --- "changes" summary cells, merged, styled like Excel's "Good" cell style ---
$ws.Range("E5").Formula = "=SUM(C5:C7) - SUM(D5:D7)"
$ws.Range("E5:E7").Merge()
$ws.Range("E5:E7").Style = "Good"
$ws.Range("E5:E7").HorizontalAlignment = -4108

$ws.Range("G5").Formula = "=SUM(D5:D9) - SUM(F5:F9)"
$ws.Range("G5:G9").Merge()
$ws.Range("G5:G9").Style = "Good"
$ws.Range("G5:G9").HorizontalAlignment = -4108

# --- Totals row 12 / ratio row 13 for the new column ---
$ws.Range("F12").Formula = "=SUM(F3:F11)"
$ws.Range("F13").Formula = "=B12/F12"

# --- Apply the "Calculation" style to the totals/ratio block (now spans A:F) ---
$ws.Range("A12:F13").Style = "Calculation"

# --- Column widths to accommodate the new columns ---
$ws.Columns("D:E").ColumnWidth = 25.71
$ws.Columns("F").ColumnWidth = 24.71
$ws.Columns("G").ColumnWidth = 12.57

# --- Selection as left by the editing author ---
$ws.Range("G15").Select()
